# Applies the "fixed export and fixing maps" edit:
#  - Rename the (only) worksheet from "1" to "აბაშა"
#  - Clear the census-note text in A2 (its shared string is dropped entirely)
#  - Delete the now-empty spacer row (old row 3)
#  - Delete the "1989" and "2002" data columns (old columns B and C),
#    keeping only the "2014" column, which slides left into column B

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab to the municipality name.
$ws.Name = "აბაშა"

# The census-note row (A2) becomes blank.
$ws.Range("A2").ClearContents()

# Remove the empty spacer row (row 3) entirely, shifting rows 4-6 up.
$ws.Rows("3").Delete()

# Remove the "1989" and "2002" columns (old B:C), shifting the "2014"
# column (old D) left into column B.
$ws.Range("B:C").EntireColumn.Delete()
